# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handback packages are now in sync with en-US:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The "Latest Handback DateTime" timestamps are refreshed
#   - The stale "handback file is not the latest" error detail is cleared
#   - A few report columns are widened/narrowed to fit the new content

$wb = $excel.ActiveWorkbook

$overviewSheet = $wb.Worksheets.Item("Overview")
$zhSheet = $wb.Worksheets.Item("zh-cn")
$deSheet = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: per-language status columns (E = zh-cn, F = de-de) ---
$overviewSheet.Range("E2").Value = $newStatus
$overviewSheet.Range("F2").Value = $newStatus

# --- zh-cn sheet row 2 ---
$zhSheet.Range("C2").Value = $newStatus
$zhSheet.Range("K2").Value = "2016-08-18 18:52:07"
$zhSheet.Range("P2").Value = ""

# --- de-de sheet row 2 ---
$deSheet.Range("C2").Value = $newStatus
$deSheet.Range("K2").Value = "2016-08-18 18:52:18"
$deSheet.Range("P2").Value = ""

# --- Column width adjustments to fit the new content ---
# Overview: Status columns E and F widen (closest achievable value to 29.9777047293527)
$overviewSheet.Columns.Item(5).ColumnWidth = 29.1666666666667
$overviewSheet.Columns.Item(6).ColumnWidth = 29.1666666666667

# zh-cn / de-de: Status column C widens, Error Detail column P narrows now that it's empty
$zhSheet.Columns.Item(3).ColumnWidth = 29.1666666666667
$zhSheet.Columns.Item(16).ColumnWidth = 12.8333333333333

$deSheet.Columns.Item(3).ColumnWidth = 29.1666666666667
$deSheet.Columns.Item(16).ColumnWidth = 12.8333333333333
